# "Hjemme passive tweaks lichtwark deleted values"
#
# The two "Lichtwark" trials previously held in columns O/R (subject CON)
# and AN/AQ (subject STR) replace the old "Hjemme" values that used to
# live in columns B/C/D/E, for every data row (header row 1 plus the two
# data rows). After the copy, the visible selection is shrunk down to the
# now-relevant B1:E3 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @{ Src = "O";  Dst = "B" },
    @{ Src = "R";  Dst = "C" },
    @{ Src = "AN"; Dst = "D" },
    @{ Src = "AQ"; Dst = "E" }
)

foreach ($pair in $pairs) {
    for ($row = 1; $row -le 3; $row++) {
        $srcCell = $ws.Range($pair.Src + $row)
        $dstCell = $ws.Range($pair.Dst + $row)
        $dstCell.Value2 = $srcCell.Value2
    }
}

# Shrink the saved selection from B1:AY3 down to B1:E3.
$ws.Range("B1:E3").Select()
